$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "wilds" row (row 4): rename detail text from "Selvageria" to "Terras Selvagens" ---
$ws.Range("C4").Value = "Terras Selvagens"

# --- New column F ("detail") descriptions for each faction row ---
# Added in row order (F4..F14) first, with F3 ("home" row) added last so
# the shared-string table ends up in the same append order as the source edit.
$ws.Range("F4").Value  = "Uma região inexplorada, onde a natureza reina e o perigo espreita."
$ws.Range("F5").Value  = "Lar dos guerreiros que vivem pela espada e pelo combate."
$ws.Range("F6").Value  = "Um refúgio sombrio para os que vivem nas sombras."
$ws.Range("F7").Value  = "Centro arcano de estudo e poder mágico."
$ws.Range("F8").Value  = "Onde riqueza e influência valem mais que espadas."
$ws.Range("F9").Value  = "Uma vila pacata cercada por mistérios antigos."
$ws.Range("F10").Value = "Cidade ancestral marcada por segredos e magia esquecida."
$ws.Range("F11").Value = "Terras frias onde o vento carrega histórias de guerra."
$ws.Range("F12").Value = "A cúpula que governa com sabedoria... E interesses ocultos."
$ws.Range("F13").Value = "Coração do comércio e das intrigas políticas."
$ws.Range("F14").Value = "Uma floresta viva, berço dos espíritos e da antiga magia."
$ws.Range("F3").Value  = "Sua base."

# --- Wrap the two longer descriptions that got a dedicated wrap-text style ---
$ws.Range("F4").WrapText = $true
$ws.Range("F8").WrapText = $true

# --- Row 4 grows tall enough to show the wrapped text ---
$ws.Rows(4).RowHeight = 75

# --- Column F is widened to fit the new description text ---
$ws.Columns("F").ColumnWidth = 54.42

# --- Selection moves from G6 to the newly populated F3 ---
$null = $ws.Range("F3").Select()
